$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 732.125
$ws.Range("I5").Value = 105.6
$ws.Range("K5").Value = 105.6
$ws.Range("M5").Value = 9.400000000000006

$ws.Range("H9").Value = 975.2
$ws.Range("I9").Value = 1084.091
$ws.Range("K9").Value = 1084.091
$ws.Range("M9").Value = -915.0909999999999

$ws.Range("H15").Value = 1323.2963
$ws.Range("I15").Value = 1323.2963
$ws.Range("K15").Value = 3969.8889
$ws.Range("M15").Value = -3800.8889

$ws.Range("H55").Value = 700
$ws.Range("I55").Value = 433.33334
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 433.33334
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = -219.33334
$ws.Range("N55").Value = -1928

$ws.Range("H69").Value = 35714.93
$ws.Range("J69").Value = 23000
$ws.Range("L69").Value = 69000
$ws.Range("N69").Value = -70748

$ws.Range("H72").Value = 35714.93
$ws.Range("J72").Value = 23000
$ws.Range("L72").Value = 207000
$ws.Range("N72").Value = -215736

$ws.Range("H99").Value = 290.83334
$ws.Range("I99").Value = 238
$ws.Range("J99").Value = 555
$ws.Range("K99").Value = 714
$ws.Range("L99").Value = 1665
$ws.Range("M99").Value = 784
$ws.Range("N99").Value = -4661

$ws.Range("H101").Value = 11113494
$ws.Range("I101").Value = 25004474
$ws.Range("K101").Value = 75013422
$ws.Range("M101").Value = -75011800

$ws.Range("H113").Value = 3699.875
$ws.Range("I113").Value = 3200
$ws.Range("K113").Value = 3200
$ws.Range("M113").Value = 54

$ws.Range("H125").Value = 187512320
$ws.Range("I125").Value = 333334660
$ws.Range("K125").Value = 3000011940
$ws.Range("M125").Value = -3000009480

$ws.Range("H135").Value = 357.92856
$ws.Range("I135").Value = 357.92856
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3221.35704
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -686.3570399999999
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 158.14285
$ws.Range("I5").Value = 158.14285
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 158.14285
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -46.14285000000001
$ws.Range("N5").ClearContents()

$ws.Range("H45").Value = 5284.8335
$ws.Range("I45").Value = 7194.75
$ws.Range("K45").Value = 7194.75
$ws.Range("M45").Value = -6817.75

$ws.Range("H88").Value = 1598.0625
$ws.Range("I88").Value = 1477.2222
$ws.Range("J88").Value = 1753.4286
$ws.Range("K88").Value = 1477.2222
$ws.Range("L88").Value = 1753.4286
$ws.Range("M88").Value = -1071.2222
$ws.Range("N88").Value = -2565.4286

$ws.Range("H91").Value = 1598.0625
$ws.Range("I91").Value = 1477.2222
$ws.Range("J91").Value = 1753.4286
$ws.Range("K91").Value = 1477.2222
$ws.Range("L91").Value = 1753.4286
$ws.Range("M91").Value = -73.22219999999993
$ws.Range("N91").Value = -4561.4286

$ws.Range("H122").Value = 1244.8438
$ws.Range("I122").Value = 1118.4445
$ws.Range("J122").Value = 1927.4
$ws.Range("K122").Value = 3355.3335
$ws.Range("L122").Value = 5782.200000000001
$ws.Range("M122").Value = -905.3335000000002
$ws.Range("N122").Value = -10682.2

$ws.Range("H132").Value = 2359.4443
$ws.Range("I132").Value = 2452.625
$ws.Range("K132").Value = 7357.875
$ws.Range("M132").Value = -4827.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 158.14285
$ws.Range("I4").Value = 158.14285
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 158.14285
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -43.14285000000001
$ws.Range("N4").ClearContents()

$ws.Range("H22").Value = 384.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H134").Value = 2510.9167
$ws.Range("I134").Value = 2510.9167
$ws.Range("K134").Value = 7532.750100000001
$ws.Range("M134").Value = -4997.750100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 494.33334
$ws.Range("I22").Value = 494.33334
$ws.Range("K22").Value = 494.33334
$ws.Range("M22").Value = -144.33334

$ws.Range("H31").Value = 2050
$ws.Range("I31").Value = 1100
$ws.Range("K31").Value = 1100
$ws.Range("M31").Value = -805

$ws.Range("H34").Value = 2050
$ws.Range("I34").Value = 1100
$ws.Range("K34").Value = 1100
$ws.Range("M34").Value = -898

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H58").Value = 2272.7778
$ws.Range("I58").Value = 2001.5714
$ws.Range("K58").Value = 2001.5714
$ws.Range("M58").Value = -1798.5714

$ws.Range("H62").Value = 4970.5713
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4970.5713
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4970.5713
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6218.5713

$ws.Range("H65").Value = 4970.5713
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4970.5713
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24852.8565
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31092.8565

$ws.Range("H94").Value = 88766.62
$ws.Range("I94").Value = 161351.14
$ws.Range("J94").Value = 4084.6667
$ws.Range("K94").Value = 161351.14
$ws.Range("L94").Value = 4084.6667
$ws.Range("M94").Value = -160900.14
$ws.Range("N94").Value = -4986.6667

$ws.Range("H122").Value = 2499.2856
$ws.Range("I122").Value = 1999.1666
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 5997.4998
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -3547.4998
$ws.Range("N122").Value = -21400

$ws.Range("H136").Value = 2272.7778
$ws.Range("I136").Value = 2001.5714
$ws.Range("K136").Value = 6004.7142
$ws.Range("M136").Value = -3454.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1687.8
$ws.Range("I38").Value = 1874.4445
$ws.Range("K38").Value = 5623.333500000001
$ws.Range("M38").Value = -5276.333500000001

$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744

$ws.Range("H132").Value = 5182.1665
$ws.Range("J132").Value = 7599
$ws.Range("L132").Value = 68391
$ws.Range("N132").Value = -73451

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 248.8
$ws.Range("I113").Value = 248.8
$ws.Range("K113").Value = 248.8
$ws.Range("M113").Value = 1921.2

$ws.Range("H132").Value = 3301.8
$ws.Range("I132").Value = 3327.25
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 9981.75
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -7451.75
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1843.5834
$ws.Range("I22").Value = 1393.5
$ws.Range("K22").Value = 1393.5
$ws.Range("M22").Value = -1098.5

$ws.Range("H27").Value = 1843.5834
$ws.Range("I27").Value = 1393.5
$ws.Range("K27").Value = 1393.5
$ws.Range("M27").Value = -1286.5

$ws.Range("H46").Value = 2993.4285
$ws.Range("I46").Value = 2492.3333
$ws.Range("K46").Value = 2492.3333
$ws.Range("M46").Value = -2304.3333

$ws.Range("H122").Value = 5899.185
$ws.Range("I122").Value = 4672.769
$ws.Range("J122").Value = 7038
$ws.Range("K122").Value = 14018.307
$ws.Range("L122").Value = 21114
$ws.Range("M122").Value = -11568.307
$ws.Range("N122").Value = -26014

$ws.Range("H132").Value = 7898.091
$ws.Range("I132").Value = 7898.091
$ws.Range("K132").Value = 23694.273
$ws.Range("M132").Value = -21164.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2646.9
$ws.Range("I122").Value = 2737.5
$ws.Range("J122").Value = 2284.5
$ws.Range("K122").Value = 8212.5
$ws.Range("L122").Value = 6853.5
$ws.Range("M122").Value = -5762.5
$ws.Range("N122").Value = -11753.5

$ws.Range("H132").Value = 2160
$ws.Range("I132").Value = 2192.1
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 6576.299999999999
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -4046.299999999999
$ws.Range("N132").Value = -11058.5
